$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.472.00"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "1.805.93"
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "228.47"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").Value = "0.583"
$ws.Range("E6").Value = "  +4.40%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "34.88"
$ws.Range("E8").Value = "  +5.79%  "

$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D12").Value = "2.066.01"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").Value = "11.25"
$ws.Range("E13").Value = "  +1.65%  "

$ws.Range("D14").Value = "1.804.06"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("E15").Value = "  +1.11%  "

$ws.Range("D16").Value = "34.454.44"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("E17").Value = "  +1.85%  "

$ws.Range("D18").Value = "69.13"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "0.0₃0800"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").Value = "245.98"
$ws.Range("E20").Value = "  -0.97%  "

$ws.Range("D21").Value = "11.53"
$ws.Range("E21").Value = "  +2.27%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "174.47"
$ws.Range("E24").Value = "  +5.43%  "

$ws.Range("E25").Value = "  +2.29%  "

$ws.Range("D26").Value = "7.77"
$ws.Range("E26").Value = "  +6.47%  "

$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("D28").Value = "0.120"
$ws.Range("E28").Value = "  +2.67%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "4.02"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").Value = "1.397.71"
$ws.Range("E35").Value = "  -2.02%  "

$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("E37").Value = "  -2.09%  "

$ws.Range("D39").Value = "0.0190"
$ws.Range("E39").Value = "  -1.26%  "

$ws.Range("D40").Value = "83.37"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("E41").Value = "  +2.91%  "

$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").Value = "13.48"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("E45").Value = "  +3.51%  "

$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("D47").Value = "5.98"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("D48").Value = "1.966.23"
$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("D49").Value = "105.00"
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("E51").Value = "  -0.12%  "
